$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.640.61'
$ws.Cells.Item(2, 5).Value = '  +0.99%  '
$ws.Cells.Item(3, 4).Value = '1.853.73'
$ws.Cells.Item(3, 5).Value = '  +0.52%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '264.84'
$ws.Cells.Item(5, 5).Value = '  +1.71%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5269'
$ws.Cells.Item(7, 5).Value = '  +0.35%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3254'
$ws.Cells.Item(8, 5).Value = '  +0.51%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06798'
$ws.Cells.Item(9, 5).Value = '  +0.72%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '19.05'
$ws.Cells.Item(10, 5).Value = '  +0.73%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.7832'
$ws.Cells.Item(11, 5).Value = '  +1.38%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.07796'
$ws.Cells.Item(12, 5).Value = '  +1.38%  '
$ws.Cells.Item(13, 4).Value = '1.859.78'
$ws.Cells.Item(13, 5).Value = '  +0.56%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '88.63'
$ws.Cells.Item(14, 5).Value = '  -0.44%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '5.029'
$ws.Cells.Item(15, 5).Value = '  -0.10%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.003'
$ws.Cells.Item(16, 5).Value = '  +0.15%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '14.03'
$ws.Cells.Item(17, 5).Value = '  -0.87%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.000007976'
$ws.Cells.Item(18, 5).Value = '  +1.18%  '
$ws.Cells.Item(19, 5).Value = '  -0.03%  '
$ws.Cells.Item(20, 4).Value = '26.654.71'
$ws.Cells.Item(20, 5).Value = '  +0.97%  '
$ws.Cells.Item(21, 5).Value = '  +2.36%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '9.514'
$ws.Cells.Item(22, 5).Value = '  +0.38%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '6.010'
$ws.Cells.Item(23, 5).Value = '  +1.36%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '142.92'
$ws.Cells.Item(24, 5).Value = '  -0.94%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.184'
$ws.Cells.Item(25, 5).Value = '  -6.53%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.684'
$ws.Cells.Item(26, 5).Value = '  +2.33%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '17.06'
$ws.Cells.Item(27, 5).Value = '  +0.71%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '112.12'
$ws.Cells.Item(28, 5).Value = '  +0.53%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '4.202'
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.119'
$ws.Cells.Item(30, 5).Value = '  -0.44%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.08733'
$ws.Cells.Item(31, 5).Value = '  -0.77%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.04849'
$ws.Cells.Item(32, 5).Value = '  +0.05%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.7232'
$ws.Cells.Item(33, 5).Value = '  +6.06%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.132'
$ws.Cells.Item(34, 5).Value = '  -0.13%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.881'
$ws.Cells.Item(35, 5).Value = '  +0.99%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '3.121'
$ws.Cells.Item(36, 5).Value = '  +0.27%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.270'
$ws.Cells.Item(37, 5).Value = '  +2.53%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.01797'
$ws.Cells.Item(38, 5).Value = '  +0.14%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.4885'
$ws.Cells.Item(39, 5).Value = '  -0.74%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.9055'
$ws.Cells.Item(40, 5).Value = '  +0.59%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '110.96'
$ws.Cells.Item(41, 5).Value = '  -1.44%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '5.982'
$ws.Cells.Item(42, 5).Value = '  -3.08%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.001'
$ws.Cells.Item(43, 5).Value = '  -0.03%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '7.696'
$ws.Cells.Item(44, 5).Value = '  -0.86%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.4210'
$ws.Cells.Item(45, 5).Value = '  +0.37%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.05892'
$ws.Cells.Item(46, 5).Value = '  +0.18%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '9.066'
$ws.Cells.Item(47, 5).Value = '  -0.36%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.1238'
$ws.Cells.Item(48, 5).Value = '  -1.95%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '35.12'
$ws.Cells.Item(49, 5).Value = '  -0.79%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.8900'
$ws.Cells.Item(50, 5).Value = '  +3.44%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '60.09'
$ws.Cells.Item(51, 5).Value = '  +1.21%  '
